# dumpLog.xlsx update: add a new dump log entry row (Dump20160210) and
# move the active selection to C7, matching a new DB dump log entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dump log entry on row 6 (FILENAME / DESCRIPTION columns B / C)
$ws.Range("B6").Value = "Dump20160210"
$ws.Range("C6").Value = "Added new tables (MASTER_BRANCH, MASTER_MODULE, USER_MANAGEMENT_ACCESS)"

# The description is long, so wrap the text in the description cell
$ws.Range("C6").WrapText = $true

# Move / leave the active selection at C7, as recorded in the sheet view
$ws.Range("C7").Select()
